$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.587.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.34%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.115.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.27%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "530.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.81%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.114.69"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.29%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.470"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.65%  "

$ws.Range("E11").Value = "  +0.67%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.410"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.48%  "

$ws.Range("E13").Value = "  +1.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.648.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.76%  "

$ws.Range("E16").Value = "  +1.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "57.698.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.39%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.115.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.11%  "

$ws.Range("E19").Value = "  +2.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.50%  "

$ws.Range("E21").Value = "  +3.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "360.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.16%  "

$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("E24").Value = "  +2.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.504"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.86%  "

$ws.Range("E26").Value = "  -0.17%  "

$ws.Range("E27").Value = "  -0.12%  "

$ws.Range("E28").Value = "  -3.64%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.59%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.28%  "

$ws.Range("E33").Value = "  +3.90%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "159.21"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.66%  "

$ws.Range("E35").Value = "  -0.98%  "

$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("E37").Value = "  +4.46%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.48"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.91%  "

$ws.Range("E39").Value = "  +3.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0668"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.488.81"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.81%  "

$ws.Range("E42").Value = "  -3.42%  "

$ws.Range("E43").Value = "  -0.49%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "37.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.24%  "

$ws.Range("E45").Value = "  +1.53%  "

$ws.Range("E46").Value = "  -0.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.977"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.87%  "

$ws.Range("E48").Value = "  +1.43%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.64"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.737"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.58%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0911"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.21%  "
